$wb = $excel.ActiveWorkbook

$wb.Worksheets("ALC").Range("H9").Value = 388.05
$wb.Worksheets("ALC").Range("I9").Value = 313.93332
$wb.Worksheets("ALC").Range("J9").Value = 610.4
$wb.Worksheets("ALC").Range("K9").Value = 313.93332
$wb.Worksheets("ALC").Range("L9").Value = 610.4
$wb.Worksheets("ALC").Range("M9").Value = -144.93332
$wb.Worksheets("ALC").Range("N9").Value = -948.4
$wb.Worksheets("ALC").Range("H12").Value = 1674.75
$wb.Worksheets("ALC").Range("I12").Value = 1567
$wb.Worksheets("ALC").Range("J12").Value = 1998
$wb.Worksheets("ALC").Range("K12").Value = 1567
$wb.Worksheets("ALC").Range("L12").Value = 1998
$wb.Worksheets("ALC").Range("M12").Value = -1397
$wb.Worksheets("ALC").Range("N12").Value = -2338
$wb.Worksheets("ALC").Range("H15").Value = 21300888
$wb.Worksheets("ALC").Range("I15").Value = 21300888
$wb.Worksheets("ALC").Range("K15").Value = 63902664
$wb.Worksheets("ALC").Range("M15").Value = -63902495
$wb.Worksheets("ALC").Range("H19").Value = 1010.94446
$wb.Worksheets("ALC").Range("I19").Value = 425.81818
$wb.Worksheets("ALC").Range("K19").Value = 425.81818
$wb.Worksheets("ALC").Range("M19").Value = -250.81818
$wb.Worksheets("ALC").Range("H21").Value = 27640.68
$wb.Worksheets("ALC").Range("I21").Value = 25254.25
$wb.Worksheets("ALC").Range("K21").Value = 25254.25
$wb.Worksheets("ALC").Range("M21").Value = -24786.25
$wb.Worksheets("ALC").Range("H23").Value = 27640.68
$wb.Worksheets("ALC").Range("I23").Value = 25254.25
$wb.Worksheets("ALC").Range("K23").Value = 25254.25
$wb.Worksheets("ALC").Range("M23").Value = -25020.25
$wb.Worksheets("ALC").Range("H33").Value = 51.142857
$wb.Worksheets("ALC").Range("I33").Value = 52.333332
$wb.Worksheets("ALC").Range("K33").Value = 52.333332
$wb.Worksheets("ALC").Range("M33").Value = 176.666668
$wb.Worksheets("ALC").Range("H38").Value = 29
$wb.Worksheets("ALC").Range("I38").Value = 29
$wb.Worksheets("ALC").Range("K38").Value = 87
$wb.Worksheets("ALC").Range("M38").Value = 285
$wb.Worksheets("ALC").Range("H43").Value = 7799.8
$wb.Worksheets("ALC").Range("I43").Value = 0
$wb.Worksheets("ALC").Range("J43").Value = 7799.8
$wb.Worksheets("ALC").Range("K43").Value = 0
$wb.Worksheets("ALC").Range("L43").Value = 7799.8
$wb.Worksheets("ALC").Range("M43").ClearContents()
$wb.Worksheets("ALC").Range("N43").Value = -7937.8
$wb.Worksheets("ALC").Range("H58").Value = 45456610
$wb.Worksheets("ALC").Range("I58").Value = 183.57143
$wb.Worksheets("ALC").Range("J58").Value = 125005350
$wb.Worksheets("ALC").Range("K58").Value = 550.71429
$wb.Worksheets("ALC").Range("L58").Value = 375016050
$wb.Worksheets("ALC").Range("M58").Value = -400.71429
$wb.Worksheets("ALC").Range("N58").Value = -375016350
$wb.Worksheets("ALC").Range("H62").Value = 1300.4
$wb.Worksheets("ALC").Range("I62").Value = 1300.4
$wb.Worksheets("ALC").Range("J62").Value = 0
$wb.Worksheets("ALC").Range("K62").Value = 1300.4
$wb.Worksheets("ALC").Range("L62").Value = 0
$wb.Worksheets("ALC").Range("M62").Value = -676.4000000000001
$wb.Worksheets("ALC").Range("N62").ClearContents()
$wb.Worksheets("ALC").Range("H65").Value = 1300.4
$wb.Worksheets("ALC").Range("I65").Value = 1300.4
$wb.Worksheets("ALC").Range("J65").Value = 0
$wb.Worksheets("ALC").Range("K65").Value = 6502
$wb.Worksheets("ALC").Range("L65").Value = 0
$wb.Worksheets("ALC").Range("M65").Value = -3382
$wb.Worksheets("ALC").Range("N65").ClearContents()
$wb.Worksheets("ALC").Range("H92").Value = 45455308
$wb.Worksheets("ALC").Range("I92").Value = 587.95
$wb.Worksheets("ALC").Range("J92").Value = 500002500
$wb.Worksheets("ALC").Range("K92").Value = 587.95
$wb.Worksheets("ALC").Range("L92").Value = 500002500
$wb.Worksheets("ALC").Range("M92").Value = 660.05
$wb.Worksheets("ALC").Range("N92").Value = -500004996
$wb.Worksheets("ALC").Range("H107").Value = 39062990
$wb.Worksheets("ALC").Range("I107").Value = 9615930
$wb.Worksheets("ALC").Range("J107").Value = 166666930
$wb.Worksheets("ALC").Range("K107").Value = 9615930
$wb.Worksheets("ALC").Range("L107").Value = 166666930
$wb.Worksheets("ALC").Range("M107").Value = -9614010
$wb.Worksheets("ALC").Range("N107").Value = -166670770
$wb.Worksheets("ALC").Range("H121").Value = 450
$wb.Worksheets("ALC").Range("J121").Value = 0
$wb.Worksheets("ALC").Range("L121").Value = 0
$wb.Worksheets("ALC").Range("N121").ClearContents()
$wb.Worksheets("ALC").Range("H133").Value = 0
$wb.Worksheets("ALC").Range("J133").Value = 0
$wb.Worksheets("ALC").Range("L133").Value = 0
$wb.Worksheets("ALC").Range("N133").ClearContents()
$wb.Worksheets("ALC").Range("H137").Value = 5417.6855
$wb.Worksheets("ALC").Range("I137").Value = 2713.2778
$wb.Worksheets("ALC").Range("K137").Value = 8139.8334
$wb.Worksheets("ALC").Range("M137").Value = -5589.8334
$wb.Worksheets("ALC").Range("H138").Value = 1067068.9
$wb.Worksheets("ALC").Range("I138").Value = 1551.9791
$wb.Worksheets("ALC").Range("J138").Value = 2178912.8
$wb.Worksheets("ALC").Range("K138").Value = 4655.9373
$wb.Worksheets("ALC").Range("L138").Value = 6536738.399999999
$wb.Worksheets("ALC").Range("M138").Value = 484.0627000000004
$wb.Worksheets("ALC").Range("N138").Value = -6547018.399999999
$wb.Worksheets("ARM").Range("H2").Value = 19608878
$wb.Worksheets("ARM").Range("I2").Value = 853.6842
$wb.Worksheets("ARM").Range("K2").Value = 853.6842
$wb.Worksheets("ARM").Range("M2").Value = -740.6842
$wb.Worksheets("ARM").Range("H32").Value = 2252580.2
$wb.Worksheets("ARM").Range("I32").Value = 2385197.8
$wb.Worksheets("ARM").Range("K32").Value = 2385197.8
$wb.Worksheets("ARM").Range("M32").Value = -2384910.8
$wb.Worksheets("ARM").Range("H45").Value = 5128.1177
$wb.Worksheets("ARM").Range("J45").Value = 8003.5557
$wb.Worksheets("ARM").Range("L45").Value = 8003.5557
$wb.Worksheets("ARM").Range("N45").Value = -8757.555700000001
$wb.Worksheets("ARM").Range("H61").Value = 34489244
$wb.Worksheets("ARM").Range("I61").Value = 3131.95
$wb.Worksheets("ARM").Range("K61").Value = 3131.95
$wb.Worksheets("ARM").Range("M61").Value = -2919.95
$wb.Worksheets("ARM").Range("H74").Value = 4025.6365
$wb.Worksheets("ARM").Range("I74").Value = 2667.2
$wb.Worksheets("ARM").Range("K74").Value = 2667.2
$wb.Worksheets("ARM").Range("M74").Value = -1793.2
$wb.Worksheets("ARM").Range("H77").Value = 4025.6365
$wb.Worksheets("ARM").Range("I77").Value = 2667.2
$wb.Worksheets("ARM").Range("K77").Value = 13336
$wb.Worksheets("ARM").Range("M77").Value = -8968
$wb.Worksheets("ARM").Range("H97").Value = 8335340.5
$wb.Worksheets("ARM").Range("I97").Value = 1652.5714
$wb.Worksheets("ARM").Range("J97").Value = 27780612
$wb.Worksheets("ARM").Range("K97").Value = 1652.5714
$wb.Worksheets("ARM").Range("L97").Value = 27780612
$wb.Worksheets("ARM").Range("M97").Value = -1156.5714
$wb.Worksheets("ARM").Range("N97").Value = -27781604
$wb.Worksheets("ARM").Range("H110").Value = 23810882
$wb.Worksheets("ARM").Range("I110").Value = 1311.25
$wb.Worksheets("ARM").Range("K110").Value = 1311.25
$wb.Worksheets("ARM").Range("M110").Value = 733.75
$wb.Worksheets("ARM").Range("H116").Value = 19608878
$wb.Worksheets("ARM").Range("I116").Value = 853.6842
$wb.Worksheets("ARM").Range("K116").Value = 853.6842
$wb.Worksheets("ARM").Range("M116").Value = 1440.3158
$wb.Worksheets("ARM").Range("H122").Value = 2490.037
$wb.Worksheets("ARM").Range("I122").Value = 2015.579
$wb.Worksheets("ARM").Range("J122").Value = 3616.875
$wb.Worksheets("ARM").Range("K122").Value = 6046.737
$wb.Worksheets("ARM").Range("L122").Value = 10850.625
$wb.Worksheets("ARM").Range("M122").Value = -3596.737
$wb.Worksheets("ARM").Range("N122").Value = -15750.625
$wb.Worksheets("ARM").Range("H123").Value = 48224.223
$wb.Worksheets("ARM").Range("J123").Value = 48224.223
$wb.Worksheets("ARM").Range("L123").Value = 48224.223
$wb.Worksheets("ARM").Range("N123").Value = -58024.223
$wb.Worksheets("ARM").Range("H132").Value = 4776.3774
$wb.Worksheets("ARM").Range("I132").Value = 2246.9355
$wb.Worksheets("ARM").Range("J132").Value = 8340.591
$wb.Worksheets("ARM").Range("K132").Value = 6740.806500000001
$wb.Worksheets("ARM").Range("L132").Value = 25021.773
$wb.Worksheets("ARM").Range("M132").Value = -4210.806500000001
$wb.Worksheets("ARM").Range("N132").Value = -30081.773
$wb.Worksheets("ARM").Range("H136").Value = 34489244
$wb.Worksheets("ARM").Range("I136").Value = 3131.95
$wb.Worksheets("ARM").Range("K136").Value = 9395.849999999999
$wb.Worksheets("ARM").Range("M136").Value = -6845.849999999999
$wb.Worksheets("BSM").Range("H3").Value = 19608878
$wb.Worksheets("BSM").Range("I3").Value = 853.6842
$wb.Worksheets("BSM").Range("K3").Value = 853.6842
$wb.Worksheets("BSM").Range("M3").Value = -739.6842
$wb.Worksheets("BSM").Range("H86").Value = 7355067
$wb.Worksheets("BSM").Range("I86").Value = 10418446
$wb.Worksheets("BSM").Range("K86").Value = 10418446
$wb.Worksheets("BSM").Range("M86").Value = -10417323
$wb.Worksheets("BSM").Range("H89").Value = 7355067
$wb.Worksheets("BSM").Range("I89").Value = 10418446
$wb.Worksheets("BSM").Range("K89").Value = 52092230
$wb.Worksheets("BSM").Range("M89").Value = -52086614
$wb.Worksheets("BSM").Range("H94").Value = 2238.1333
$wb.Worksheets("BSM").Range("I94").Value = 1668.5
$wb.Worksheets("BSM").Range("J94").Value = 4516.6665
$wb.Worksheets("BSM").Range("K94").Value = 1668.5
$wb.Worksheets("BSM").Range("L94").Value = 4516.6665
$wb.Worksheets("BSM").Range("M94").Value = -1217.5
$wb.Worksheets("BSM").Range("N94").Value = -5418.6665
$wb.Worksheets("BSM").Range("H134").Value = 7582747.5
$wb.Worksheets("BSM").Range("I134").Value = 17859594
$wb.Worksheets("BSM").Range("K134").Value = 53578782
$wb.Worksheets("BSM").Range("M134").Value = -53576247
$wb.Worksheets("CRP").Range("H16").Value = 6175.6875
$wb.Worksheets("CRP").Range("I16").Value = 3994.75
$wb.Worksheets("CRP").Range("J16").Value = 6902.6665
$wb.Worksheets("CRP").Range("K16").Value = 3994.75
$wb.Worksheets("CRP").Range("L16").Value = 6902.6665
$wb.Worksheets("CRP").Range("M16").Value = -3707.75
$wb.Worksheets("CRP").Range("N16").Value = -7476.6665
$wb.Worksheets("CRP").Range("H52").Value = 97036.336
$wb.Worksheets("CRP").Range("I52").Value = 80000
$wb.Worksheets("CRP").Range("J52").Value = 105554.5
$wb.Worksheets("CRP").Range("K52").Value = 80000
$wb.Worksheets("CRP").Range("L52").Value = 105554.5
$wb.Worksheets("CRP").Range("M52").Value = -79706
$wb.Worksheets("CRP").Range("N52").Value = -106142.5
$wb.Worksheets("CRP").Range("H62").Value = 6949611
$wb.Worksheets("CRP").Range("I62").Value = 15630373
$wb.Worksheets("CRP").Range("J62").Value = 5001.1
$wb.Worksheets("CRP").Range("K62").Value = 15630373
$wb.Worksheets("CRP").Range("L62").Value = 5001.1
$wb.Worksheets("CRP").Range("M62").Value = -15629749
$wb.Worksheets("CRP").Range("N62").Value = -6249.1
$wb.Worksheets("CRP").Range("H65").Value = 6949611
$wb.Worksheets("CRP").Range("I65").Value = 15630373
$wb.Worksheets("CRP").Range("J65").Value = 5001.1
$wb.Worksheets("CRP").Range("K65").Value = 78151865
$wb.Worksheets("CRP").Range("L65").Value = 25005.5
$wb.Worksheets("CRP").Range("M65").Value = -78148745
$wb.Worksheets("CRP").Range("N65").Value = -31245.5
$wb.Worksheets("CRP").Range("H70").Value = 0
$wb.Worksheets("CRP").Range("J70").Value = 0
$wb.Worksheets("CRP").Range("L70").Value = 0
$wb.Worksheets("CRP").Range("N70").ClearContents()
$wb.Worksheets("CRP").Range("H73").Value = 0
$wb.Worksheets("CRP").Range("J73").Value = 0
$wb.Worksheets("CRP").Range("L73").Value = 0
$wb.Worksheets("CRP").Range("N73").ClearContents()
$wb.Worksheets("CRP").Range("H86").Value = 16452974
$wb.Worksheets("CRP").Range("I86").Value = 22327538
$wb.Worksheets("CRP").Range("J86").Value = 4193.6
$wb.Worksheets("CRP").Range("K86").Value = 22327538
$wb.Worksheets("CRP").Range("L86").Value = 4193.6
$wb.Worksheets("CRP").Range("M86").Value = -22326415
$wb.Worksheets("CRP").Range("N86").Value = -6439.6
$wb.Worksheets("CRP").Range("H89").Value = 16452974
$wb.Worksheets("CRP").Range("I89").Value = 22327538
$wb.Worksheets("CRP").Range("J89").Value = 4193.6
$wb.Worksheets("CRP").Range("K89").Value = 111637690
$wb.Worksheets("CRP").Range("L89").Value = 20968
$wb.Worksheets("CRP").Range("M89").Value = -111632074
$wb.Worksheets("CRP").Range("N89").Value = -32200
$wb.Worksheets("CRP").Range("H113").Value = 6175.6875
$wb.Worksheets("CRP").Range("I113").Value = 3994.75
$wb.Worksheets("CRP").Range("J113").Value = 6902.6665
$wb.Worksheets("CRP").Range("K113").Value = 3994.75
$wb.Worksheets("CRP").Range("L113").Value = 6902.6665
$wb.Worksheets("CRP").Range("M113").Value = -1824.75
$wb.Worksheets("CRP").Range("N113").Value = -11242.6665
$wb.Worksheets("CRP").Range("H134").Value = 7846.8184
$wb.Worksheets("CRP").Range("I134").Value = 3332.3333
$wb.Worksheets("CRP").Range("K134").Value = 9996.999899999999
$wb.Worksheets("CRP").Range("M134").Value = -7461.999899999999
$wb.Worksheets("CRP").Range("H140").Value = 30000
$wb.Worksheets("CRP").Range("I140").Value = 30000
$wb.Worksheets("CRP").Range("J140").Value = 0
$wb.Worksheets("CRP").Range("K140").Value = 30000
$wb.Worksheets("CRP").Range("L140").Value = 0
$wb.Worksheets("CRP").Range("N140").ClearContents()
$wb.Worksheets("CRP").Range("M140").Value = -24820
$wb.Worksheets("CUL").Range("H4").Value = 62156480
$wb.Worksheets("CUL").Range("I4").Value = 78509990
$wb.Worksheets("CUL").Range("J4").Value = 4919202
$wb.Worksheets("CUL").Range("K4").Value = 235529970
$wb.Worksheets("CUL").Range("L4").Value = 14757606
$wb.Worksheets("CUL").Range("M4").Value = -235529858
$wb.Worksheets("CUL").Range("N4").Value = -14757830
$wb.Worksheets("CUL").Range("H5").Value = 1504.4348
$wb.Worksheets("CUL").Range("I5").Value = 682.7646999999999
$wb.Worksheets("CUL").Range("J5").Value = 3832.5
$wb.Worksheets("CUL").Range("K5").Value = 2048.2941
$wb.Worksheets("CUL").Range("L5").Value = 11497.5
$wb.Worksheets("CUL").Range("M5").Value = -1936.2941
$wb.Worksheets("CUL").Range("N5").Value = -11721.5
$wb.Worksheets("CUL").Range("H17").Value = 4750
$wb.Worksheets("CUL").Range("J17").Value = 4750
$wb.Worksheets("CUL").Range("L17").Value = 14250
$wb.Worksheets("CUL").Range("N17").Value = -14588
$wb.Worksheets("CUL").Range("H113").Value = 7016.4287
$wb.Worksheets("CUL").Range("J113").Value = 12496.571
$wb.Worksheets("CUL").Range("L113").Value = 37489.713
$wb.Worksheets("CUL").Range("N113").Value = -41829.713
$wb.Worksheets("CUL").Range("H135").Value = 1504.4348
$wb.Worksheets("CUL").Range("I135").Value = 682.7646999999999
$wb.Worksheets("CUL").Range("J135").Value = 3832.5
$wb.Worksheets("CUL").Range("K135").Value = 6144.882299999999
$wb.Worksheets("CUL").Range("L135").Value = 34492.5
$wb.Worksheets("CUL").Range("M135").Value = -3609.882299999999
$wb.Worksheets("CUL").Range("N135").Value = -39562.5
$wb.Worksheets("GSM").Range("H10").Value = 8428.571
$wb.Worksheets("GSM").Range("J10").Value = 7600
$wb.Worksheets("GSM").Range("L10").Value = 7600
$wb.Worksheets("GSM").Range("N10").Value = -7938
$wb.Worksheets("GSM").Range("H62").Value = 70077
$wb.Worksheets("GSM").Range("I62").Value = 70077
$wb.Worksheets("GSM").Range("K62").Value = 70077
$wb.Worksheets("GSM").Range("M62").Value = -69391
$wb.Worksheets("GSM").Range("H65").Value = 70077
$wb.Worksheets("GSM").Range("I65").Value = 70077
$wb.Worksheets("GSM").Range("K65").Value = 210231
$wb.Worksheets("GSM").Range("M65").Value = -206799
$wb.Worksheets("GSM").Range("H93").Value = 49967.5
$wb.Worksheets("GSM").Range("J93").Value = 49967.5
$wb.Worksheets("GSM").Range("L93").Value = 49967.5
$wb.Worksheets("GSM").Range("N93").Value = -53711.5
$wb.Worksheets("GSM").Range("H97").Value = 1806.35
$wb.Worksheets("GSM").Range("J97").Value = 3350.5
$wb.Worksheets("GSM").Range("L97").Value = 3350.5
$wb.Worksheets("GSM").Range("N97").Value = -4342.5
$wb.Worksheets("GSM").Range("H102").Value = 3823.7334
$wb.Worksheets("GSM").Range("I102").Value = 3823.7334
$wb.Worksheets("GSM").Range("K102").Value = 3823.7334
$wb.Worksheets("GSM").Range("M102").Value = -2201.7334
$wb.Worksheets("GSM").Range("H125").Value = 85000
$wb.Worksheets("GSM").Range("J125").Value = 85000
$wb.Worksheets("GSM").Range("L125").Value = 85000
$wb.Worksheets("GSM").Range("N125").Value = -89920
$wb.Worksheets("GSM").Range("H132").Value = 5051.4443
$wb.Worksheets("GSM").Range("I132").Value = 2496
$wb.Worksheets("GSM").Range("K132").Value = 7488
$wb.Worksheets("GSM").Range("M132").Value = -4958
$wb.Worksheets("LTW").Range("H7").Value = 4856.0557
$wb.Worksheets("LTW").Range("I7").Value = 3814.9285
$wb.Worksheets("LTW").Range("K7").Value = 3814.9285
$wb.Worksheets("LTW").Range("M7").Value = -3702.9285
$wb.Worksheets("LTW").Range("H22").Value = 998.1053000000001
$wb.Worksheets("LTW").Range("I22").Value = 311.2
$wb.Worksheets("LTW").Range("K22").Value = 311.2
$wb.Worksheets("LTW").Range("M22").Value = -16.19999999999999
$wb.Worksheets("LTW").Range("H27").Value = 998.1053000000001
$wb.Worksheets("LTW").Range("I27").Value = 311.2
$wb.Worksheets("LTW").Range("K27").Value = 311.2
$wb.Worksheets("LTW").Range("M27").Value = -204.2
$wb.Worksheets("LTW").Range("H40").Value = 4286.032
$wb.Worksheets("LTW").Range("I40").Value = 3416.9092
$wb.Worksheets("LTW").Range("K40").Value = 3416.9092
$wb.Worksheets("LTW").Range("M40").Value = -3280.9092
$wb.Worksheets("LTW").Range("H46").Value = 2156576.5
$wb.Worksheets("LTW").Range("I46").Value = 2653867.8
$wb.Worksheets("LTW").Range("K46").Value = 2653867.8
$wb.Worksheets("LTW").Range("M46").Value = -2653679.8
$wb.Worksheets("LTW").Range("H68").Value = 6557.2144
$wb.Worksheets("LTW").Range("I68").Value = 4843
$wb.Worksheets("LTW").Range("J68").Value = 8271.429
$wb.Worksheets("LTW").Range("K68").Value = 4843
$wb.Worksheets("LTW").Range("L68").Value = 8271.429
$wb.Worksheets("LTW").Range("M68").Value = -4094
$wb.Worksheets("LTW").Range("N68").Value = -9769.429
$wb.Worksheets("LTW").Range("H71").Value = 6557.2144
$wb.Worksheets("LTW").Range("I71").Value = 4843
$wb.Worksheets("LTW").Range("J71").Value = 8271.429
$wb.Worksheets("LTW").Range("K71").Value = 24215
$wb.Worksheets("LTW").Range("L71").Value = 41357.145
$wb.Worksheets("LTW").Range("M71").Value = -20471
$wb.Worksheets("LTW").Range("N71").Value = -48845.145
$wb.Worksheets("LTW").Range("H93").Value = 1272.6666
$wb.Worksheets("LTW").Range("I93").Value = 909.1667
$wb.Worksheets("LTW").Range("J93").Value = 1999.6666
$wb.Worksheets("LTW").Range("K93").Value = 909.1667
$wb.Worksheets("LTW").Range("L93").Value = 1999.6666
$wb.Worksheets("LTW").Range("M93").Value = 338.8333
$wb.Worksheets("LTW").Range("N93").Value = -4495.6666
$wb.Worksheets("LTW").Range("H100").Value = 3817.875
$wb.Worksheets("LTW").Range("I100").Value = 3863.2856
$wb.Worksheets("LTW").Range("J100").Value = 3500
$wb.Worksheets("LTW").Range("K100").Value = 3863.2856
$wb.Worksheets("LTW").Range("L100").Value = 3500
$wb.Worksheets("LTW").Range("M100").Value = -3322.2856
$wb.Worksheets("LTW").Range("N100").Value = -4582
$wb.Worksheets("LTW").Range("H119").Value = 90000
$wb.Worksheets("LTW").Range("J119").Value = 90000
$wb.Worksheets("LTW").Range("L119").Value = 90000
$wb.Worksheets("LTW").Range("N119").Value = -99676
$wb.Worksheets("LTW").Range("H125").Value = 55998.332
$wb.Worksheets("LTW").Range("J125").Value = 55998.332
$wb.Worksheets("LTW").Range("L125").Value = 55998.332
$wb.Worksheets("LTW").Range("N125").Value = -65838.33199999999
$wb.Worksheets("LTW").Range("H126").Value = 4856.0557
$wb.Worksheets("LTW").Range("I126").Value = 3814.9285
$wb.Worksheets("LTW").Range("K126").Value = 11444.7855
$wb.Worksheets("LTW").Range("M126").Value = -8974.7855
$wb.Worksheets("LTW").Range("H129").Value = 95164.5
$wb.Worksheets("LTW").Range("J129").Value = 95164.5
$wb.Worksheets("LTW").Range("L129").Value = 95164.5
$wb.Worksheets("LTW").Range("N129").Value = -105164.5
$wb.Worksheets("LTW").Range("H136").Value = 11633.48
$wb.Worksheets("LTW").Range("J136").Value = 16512.129
$wb.Worksheets("LTW").Range("L136").Value = 49536.387
$wb.Worksheets("LTW").Range("N136").Value = -54636.387
$wb.Worksheets("WVR").Range("H68").Value = 41499.5
$wb.Worksheets("WVR").Range("J68").Value = 41499.5
$wb.Worksheets("WVR").Range("L68").Value = 41499.5
$wb.Worksheets("WVR").Range("N68").Value = -43121.5
$wb.Worksheets("WVR").Range("H71").Value = 41499.5
$wb.Worksheets("WVR").Range("J71").Value = 41499.5
$wb.Worksheets("WVR").Range("L71").Value = 124498.5
$wb.Worksheets("WVR").Range("N71").Value = -132610.5
$wb.Worksheets("WVR").Range("H74").Value = 5787
$wb.Worksheets("WVR").Range("J74").Value = 6418.1665
$wb.Worksheets("WVR").Range("L74").Value = 6418.1665
$wb.Worksheets("WVR").Range("N74").Value = -8290.166499999999
$wb.Worksheets("WVR").Range("H77").Value = 5787
$wb.Worksheets("WVR").Range("J77").Value = 6418.1665
$wb.Worksheets("WVR").Range("L77").Value = 19254.4995
$wb.Worksheets("WVR").Range("N77").Value = -28614.4995
$wb.Worksheets("WVR").Range("H96").Value = 1108.25
$wb.Worksheets("WVR").Range("I96").Value = 1066.3334
$wb.Worksheets("WVR").Range("J96").Value = 1234
$wb.Worksheets("WVR").Range("K96").Value = 1066.3334
$wb.Worksheets("WVR").Range("L96").Value = 1234
$wb.Worksheets("WVR").Range("M96").Value = 306.6666
$wb.Worksheets("WVR").Range("N96").Value = -3980
$wb.Worksheets("WVR").Range("H113").Value = 8220.861000000001
$wb.Worksheets("WVR").Range("I113").Value = 12867.81
$wb.Worksheets("WVR").Range("K113").Value = 38603.43
$wb.Worksheets("WVR").Range("M113").Value = -36433.43
$wb.Worksheets("WVR").Range("H119").Value = 67376.664
$wb.Worksheets("WVR").Range("J119").Value = 67376.664
$wb.Worksheets("WVR").Range("L119").Value = 67376.664
$wb.Worksheets("WVR").Range("N119").Value = -77052.664
$wb.Worksheets("WVR").Range("H128").Value = 0
$wb.Worksheets("WVR").Range("J128").Value = 0
$wb.Worksheets("WVR").Range("L128").Value = 0
$wb.Worksheets("WVR").Range("N128").ClearContents()
$wb.Worksheets("WVR").Range("H132").Value = 4701.2764
$wb.Worksheets("WVR").Range("I132").Value = 4278.657
$wb.Worksheets("WVR").Range("K132").Value = 12835.971
$wb.Worksheets("WVR").Range("M132").Value = -10305.971
$wb.Worksheets("WVR").Range("H136").Value = 20616188
$wb.Worksheets("WVR").Range("I136").Value = 45456772
$wb.Worksheets("WVR").Range("K136").Value = 136370316
$wb.Worksheets("WVR").Range("M136").Value = -136367766
